# Insert one new weekly price record for "Zapallo italiano" (Femacal de La
# Calera) as row 461, pushing the existing rows 461:489 down to 462:490.
# This mirrors how the source's row-461 insert shifted the remaining
# historical rows down by one position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 461; Excel shifts rows 461:489 -> 462:490
# and copies the formatting (incl. the date number format on column D)
# from the row above into the freshly inserted row.
$ws.Rows.Item(461).Insert()

# Populate the new row 461 with the new record's data.
$ws.Cells.Item(461, 1).Value  = 3
$ws.Cells.Item(461, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(461, 3).Value  = "Coquimbo"
$ws.Cells.Item(461, 4).Value  = 44746
$ws.Cells.Item(461, 5).Value  = 5
$ws.Cells.Item(461, 6).Value  = 100112032
$ws.Cells.Item(461, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(461, 8).Value  = "Sin especificar"
$ws.Cells.Item(461, 9).Value  = "Primera"
$ws.Cells.Item(461, 10).Value = 254
$ws.Cells.Item(461, 11).Value = 9500
$ws.Cells.Item(461, 12).Value = 10500
$ws.Cells.Item(461, 13).Value = 10002
$ws.Cells.Item(461, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(461, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(461, 16).Value = 143
$ws.Cells.Item(461, 17).Value = 70
$ws.Cells.Item(461, 18).Value = "Hortaliza"
